$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill in the missing "Neutrality" (column G) values with 0 for all data
# rows except row 149 (Switzerland), which already has a value (1).
for ($r = 2; $r -le 175; $r++) {
    if ($r -ne 149) {
        $ws.Cells.Item($r, 7).Value = 0
    }
}

# Reflect the updated selection used while correcting the column.
$ws.Range("G2:G148").Select()
